$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.587.73"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "3.074.82"
$ws.Range("E3").Value = "  +3.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "199.23"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "619.43"
$ws.Range("E6").Value = "  +3.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  +6.80%  "
$ws.Range("D10").Value = "3.075.81"
$ws.Range("E10").Value = "  +4.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.441"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.24"
$ws.Range("E13").Value = "  +6.92%  "
$ws.Range("D14").Value = "3.609.60"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.06"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").Value = "76.373.50"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000194"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("D18").Value = "3.049.94"
$ws.Range("E18").Value = "  +3.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.60"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.94"
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.34"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.47"
$ws.Range("E22").Value = "  +9.57%  "
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").Value = "3.211.23"
$ws.Range("E24").Value = "  +3.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.55"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.35"
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.88"
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000109"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.34"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "499.28"
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("E34").Value = "  +5.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.125"
$ws.Range("E36").Value = "  +13.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.75"
$ws.Range("E37").Value = "  +2.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.72"
$ws.Range("E38").Value = "  -2.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.07"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "192.81"
$ws.Range("E40").Value = "  +6.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.378"
$ws.Range("E41").Value = "  -6.17%  "
$ws.Range("E42").Value = "  -8.64%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.14"
$ws.Range("E44").Value = "  +4.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.788"
$ws.Range("E45").Value = "  +19.39%  "
$ws.Range("E46").Value = "  +5.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.23"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  +5.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.597"
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.88"
$ws.Range("E51").Value = "  -0.02%  "

Write-Output "Applied cryptos update"